$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q4" and before the
#    existing "总计" (totals) sheet.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item(1)

$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# Reuse the header styling (bold, centered, bordered) already defined on the
# "2021-Q4" sheet so no redundant style entries are introduced.
$q4Sheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1").PasteSpecial(-4122)

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Reuse the row-index styling (column A) from the "2021-Q4" sheet.
$q4Sheet.Range("A2").Copy()
$q1Sheet.Range("A2:A7").PasteSpecial(-4122)

# Columns that hold numeric-looking text (fund code / scale / position /
# ratio / market value) must be forced to Text so values such as "001487"
# or "0.4970" keep their exact original formatting instead of being
# re-interpreted as numbers.
$q1Sheet.Range("B2:B7").NumberFormat = "@"
$q1Sheet.Range("D2:G7").NumberFormat = "@"

$fundData = @(
    @(0, "001487", "宝盈优势产业灵活配置混合",     "17.02", "91.61", "2.92", "0.4970", 9),
    @(1, "005347", "诺德量化优选6个月持有期混合",   "2.60",  "93.66", "3.07", "0.0798", 5),
    @(2, "006267", "诺德量化核心灵活配置混合A",     "1.84",  "93.91", "3.07", "0.0565", 5),
    @(3, "006268", "诺德量化核心灵活配置混合C",     "0.50",  "93.91", "3.07", "0.0154", 5),
    @(4, "001303", "银华稳利灵活配置混合A",         "0.18",  "28.88", "0.90", "0.0016", 1),
    @(5, "002323", "银华稳利灵活配置混合C",         "0.12",  "28.88", "0.90", "0.0011", 1)
)

$r = 2
foreach ($row in $fundData) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: push the existing "2021-Q4" summary
#    row down to row 3 and add a new "2022-Q1" summary row in its place.
#    (Re-fetch the sheet reference now, after the insert above, since sheet
#    positions have shifted.)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$oldDate = $totalSheet.Cells.Item(2, 2).Value()
$oldCount = $totalSheet.Cells.Item(2, 3).Value()
$oldValue = $totalSheet.Cells.Item(2, 4).Value()

# Duplicate the row-index formatting (column A, style shared with the
# header row) into row 3 before writing the shifted-down data there.
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(3, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = $oldDate
$totalSheet.Cells.Item(3, 3).Value = $oldCount
$totalSheet.Cells.Item(3, 4).Value = $oldValue

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.65
